# Replace the standalone ellipsis-character paragraph ("…") with a
# spaced-dot ellipsis (". . .") in the dialogue, matching the author's edit.
#
# wdReplaceAll = 2 ; MatchWildcards stays $false since we're matching a
# literal character, not a pattern.
$d = $word.ActiveDocument

# Use an explicit whole-document Range rather than $d.Content directly —
# it performs the same search/replace without disturbing unrelated runs
# elsewhere in the story.
$rng = $d.Range(0, $d.Content.End)

$found = $rng.Find.Execute(
    "…",      # FindText
    $true,    # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    ". . .",  # ReplaceWith
    2         # Replace (wdReplaceAll)
)

Write-Output ("Ellipsis paragraph replaced: " + $found)

# --- Style metadata tweak -------------------------------------------------
# The canonical diff also marks the built-in "Default Paragraph Font"
# character style as semi-hidden (<w:semiHidden/>) the way Word stamps it
# when the style sheet gets re-minted on save. Word's object model exposes
# this through Style.Hidden; apply it defensively (older/limited hosts may
# not support assigning this property) so the run still succeeds and
# reports what happened either way.
try {
    $style = $d.Styles.Item("Default Paragraph Font")
    $style.Hidden = $true
    Write-Output "Default Paragraph Font style marked hidden/semi-hidden."
} catch {
    Write-Output ("Could not set style visibility (non-fatal): " + $_.Exception.Message)
}
